$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (changed date) column C for rows 2-5 from 2023-10-09 (45208) to 2023-10-13 (45212)
$ws.Range("C2:C5").Value = 45212
